$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "STREET 45 PVCDWV 3`" ",
    "STREET 45 PVCDWV 2`" ",
    "22 PVCDWV 3`" ",
    "STREET 22 PVCDWV 3`" ",
    "22 PVCDWV 3`" ",
    "STREET 22 PVCDWV 2`" ",
    "COUPLING PVCDWV 4`" ",
    "COUPLING PVCDWV 3`" ",
    "COUPLING PVCDWV 2`" ",
    "COUPLING PVCDWV 1 1/2 ",
    "BUSHING PVCDWV 4 x 3",
    "BUSHING PVCDWV 3 x 2",
    "BUSHING PVCDWV 2 x 1 1/2",
    "4 TECHNO CAP ",
    "3 TECHNO CAP ",
    "2 TECHNO CAP ",
    "1 1/2 TECHNO CAP ",
    "4`" HARD CAP PVC40 ",
    "3`" HARD CAP PVC40 ",
    "20' PVC FOAM PIPE 4`" ",
    "20' PVC FOAM PIPE 3`" ",
    "20' PVC FOAM PIPE 2`" ",
    "20' PVC FOAMPIPE 1 1/2`" ",
    "GALLON PVC CEMENT ",
    "GALLON PURPLE PRIMER ",
    "TEST TEE W/PLUG PVCDWV 2 ",
    "TEST TEE W/PLUG PVCDWV 3 ",
    "TEST TEE W/PLUG PVCDWV 4 ",
    "SAN TEE PVCDWV 3`" ",
    "SAN TEE PVCDWV 3 x 2 ",
    "SAN TEE PVCDWV 3 x 1 1/2  ",
    "ST SANTEE PVCDWV 3 x 11/2 ",
    "SAN TEE PVCDWV 2 x 1 1/2 ",
    "SAN TEE PVCDWV 2`" ",
    "ST SAN TEE PVCDWV 2x1 1/2 "
)

$startRow = 17
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:A62"))

$ws.Range("A51").Select()
